$d = $word.ActiveDocument

$d.Content.Find.Execute("88×12=1056", $false, $false, $false, $false, $false, $true, 1, $false, "54×30=1620", 2) | Out-Null
$d.Content.Find.Execute("53×28=1484", $false, $false, $false, $false, $false, $true, 1, $false, "55×91=5005", 2) | Out-Null
$d.Content.Find.Execute("24×23=552", $false, $false, $false, $false, $false, $true, 1, $false, "68×98=6664", 2) | Out-Null
$d.Content.Find.Execute("81×74=5994", $false, $false, $false, $false, $false, $true, 1, $false, "52×81=4212", 2) | Out-Null
$d.Content.Find.Execute("11×82=902", $false, $false, $false, $false, $false, $true, 1, $false, "43×56=2408", 2) | Out-Null
$d.Content.Find.Execute("70×41=2870", $false, $false, $false, $false, $false, $true, 1, $false, "70×91=6370", 2) | Out-Null
$d.Content.Find.Execute("39×12=468", $false, $false, $false, $false, $false, $true, 1, $false, "41×61=2501", 2) | Out-Null
$d.Content.Find.Execute("79×28=2212", $false, $false, $false, $false, $false, $true, 1, $false, "91×75=6825", 2) | Out-Null
$d.Content.Find.Execute("21×43=903", $false, $false, $false, $false, $false, $true, 1, $false, "16×50=800", 2) | Out-Null
$d.Content.Find.Execute("31×17=527", $false, $false, $false, $false, $false, $true, 1, $false, "48×93=4464", 2) | Out-Null
$d.Content.Find.Execute("72×14=1008", $false, $false, $false, $false, $false, $true, 1, $false, "93×82=7626", 2) | Out-Null
$d.Content.Find.Execute("25×17=425", $false, $false, $false, $false, $false, $true, 1, $false, "89×54=4806", 2) | Out-Null
$d.Content.Find.Execute("16×51=816", $false, $false, $false, $false, $false, $true, 1, $false, "66×74=4884", 2) | Out-Null
$d.Content.Find.Execute("94×18=1692", $false, $false, $false, $false, $false, $true, 1, $false, "45×53=2385", 2) | Out-Null
$d.Content.Find.Execute("45×43=1935", $false, $false, $false, $false, $false, $true, 1, $false, "17×56=952", 2) | Out-Null
$d.Content.Find.Execute("53×23=1219", $false, $false, $false, $false, $false, $true, 1, $false, "64×64=4096", 2) | Out-Null
$d.Content.Find.Execute("20×98=1960", $false, $false, $false, $false, $false, $true, 1, $false, "77×71=5467", 2) | Out-Null
$d.Content.Find.Execute("43×96=4128", $false, $false, $false, $false, $false, $true, 1, $false, "87×77=6699", 2) | Out-Null
$d.Content.Find.Execute("87×61=5307", $false, $false, $false, $false, $false, $true, 1, $false, "42×15=630", 2) | Out-Null
$d.Content.Find.Execute("54×71=3834", $false, $false, $false, $false, $false, $true, 1, $false, "82×45=3690", 2) | Out-Null
$d.Content.Find.Execute("35×15=525", $false, $false, $false, $false, $false, $true, 1, $false, "12×43=516", 2) | Out-Null
$d.Content.Find.Execute("30×95=2850", $false, $false, $false, $false, $false, $true, 1, $false, "18×47=846", 2) | Out-Null
$d.Content.Find.Execute("79×67=5293", $false, $false, $false, $false, $false, $true, 1, $false, "42×76=3192", 2) | Out-Null
$d.Content.Find.Execute("68×95=6460", $false, $false, $false, $false, $false, $true, 1, $false, "44×80=3520", 2) | Out-Null
$d.Content.Find.Execute("51×15=765", $false, $false, $false, $false, $false, $true, 1, $false, "40×16=640", 2) | Out-Null
